$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2 = 'HE1-4FertHighIrrigOnCvQL41xQL36'
    3 = 'HE1-4FertHighIrrigOnCvBuster'
    4 = 'HE1-4FertHighIrrigOnCvM351'
    5 = 'HE1-4FertLowIrrigOnCvQL41xQL36'
    6 = 'HE1-4FertLowIrrigOnCvBuster'
    7 = 'HE1-4FertLowIrrigOnCvM351'
    8 = 'HE1-4FertHighIrrigOffCvQL41xQL36'
    9 = 'HE1-4FertHighIrrigOffCvBuster'
    10 = 'HE1-4FertHighIrrigOffCvM351'
    11 = 'HE1-4FertLowIrrigOffCvQL41xQL36'
    12 = 'HE1-4FertLowIrrigOffCvBuster'
    13 = 'HE1-4FertLowIrrigOffCvM351'
    14 = 'HE5-8FertHighIrrigHighCvBuster'
    15 = 'HE5-8FertHighIrrigHighCvM351'
    16 = 'HE5-8FertHighIrrigHighCvCSH13R'
    17 = 'HE5-8FertLowIrrigHighCvBuster'
    18 = 'HE5-8FertLowIrrigHighCvM351'
    19 = 'HE5-8FertLowIrrigHighCvCSH13R'
    20 = 'HE5-8FertHighIrrigLowCvBuster'
    21 = 'HE5-8FertHighIrrigLowCvM351'
    22 = 'HE5-8FertHighIrrigLowCvCSH13R'
    23 = 'HE5-8FertLowIrrigLowCvBuster'
    24 = 'HE5-8FertLowIrrigLowCvM351'
    25 = 'HE5-8FertLowIrrigLowCvCSH13R'
    26 = 'BW5CvM35-1FertHighIrrigHigh'
    27 = 'BW5CvCSH13RFertLowIrrigLow'
    28 = 'BW5CvAtx623xRTX430FertLowIrrigLow'
    29 = 'BW5CvQL41xQL36FertLowIrrigLow'
    30 = 'BW5CvCSH13RFertHighIrrigHigh'
    31 = 'BW5CvAtx623xRTX430FertHighIrrigHigh'
    32 = 'BW5CvQL41xQL36FertHighIrrigHigh'
    33 = 'BW5CvM35-1FertMedIrrigLow'
    34 = 'BW5CvCSH13RFertMedIrrigLow'
    35 = 'BW5CvAtx623xRTX430FertMedIrrigLow'
    36 = 'BW5CvQL41xQL36FertMedIrrigLow'
    37 = 'BW5CvM35-1FertLowIrrigLow'
    38 = 'BW8CvM35-1FertHighIrrigOn'
    39 = 'BW8CvCSH13RFertLowIrrigOff'
    40 = 'BW8CvATX623xRTX430FertLowIrrigOff'
    41 = 'BW8CvQL41xQL36FertLowIrrigOff'
    42 = 'BW8CvCSH13RFertHighIrrigOn'
    43 = 'BW8CvAtx623xRTx430FertHighIrrigOn'
    44 = 'BW8CvQL41xQL36FertHighIrrigOn'
    45 = 'BW8CvM35-1FertMedIrrigOff'
    46 = 'BW8CvCSH13RFertMedIrrigOff'
    47 = 'BW8CvATX623xRTX430FertMedIrrigOff'
    48 = 'BW8CvQL41xQL36FertMedIrrigOff'
    49 = 'BW8CvM35-1FertLowIrrigOff'
    50 = 'LE13FertHighIrrigOnCvBuster'
    51 = 'LE13FertHighIrrigOffCvBuster'
    52 = 'LE13FertLowirrigOffCvBuster'
    53 = 'LE13FertHighIrrigOnCvM351'
    54 = 'LE13FertHighIrrigOffCvM351'
    55 = 'LE13FertLowIrrigOffCvM351'
    56 = 'LE14CvBuster'
    57 = 'LE14CvM351'
    58 = 'LE14CvQL41xQL36'
    59 = 'LE15CvBuster'
    60 = 'LE15CvM351'
    61 = 'LE15CvQL41xQL36'
    64 = 'LE17CvBuster'
    65 = 'LE17CvM351'
    66 = 'LE17CvCSH13R'
    67 = 'LE19FertOffCvBuster'
    68 = 'LE19FertLowCvBuster'
    69 = 'LE19FertMedCvBuster'
    70 = 'LE19FertHighCvBuster'
    71 = 'LE19FertOffCvCSH13R'
    72 = 'LE19FertLowCvCSH13R'
    73 = 'LE19FertMedCvCSH13R'
    74 = 'LE19FertHighCvCSH13R'
    75 = 'LE21CvCSH13RFertLow'
    76 = 'LE21CvA35xQL36FertLow'
    77 = 'LE21CvQL39xQL36FertLow'
    78 = 'LE21CvCSH13RFertMed'
    79 = 'LE21CvA35xQL36FertMed'
    80 = 'LE21CvQL39xQL36FertMed'
    81 = 'LE21CvCSH13RFertHigh'
    82 = 'LE21CvA35xQL36FertHigh'
    83 = 'LE21CvQL39xQL36FertHigh'
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 1).Value = $updates[$row]
}
